$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "cocktail_hyper_alpha" block (row 5): replace the rich-text {a/b}_alpha
#     cell with the merged, plain-text set of values, and repurpose the
#     duplicated {a/b}_gamma cell in H5 to be the {a/b}_h cell ---
$ws.Range("D5").Value = "{a/b}_alpha={0.01, 0.1, 1.0, 5.0}"
$ws.Range("H5").Value = "{a/b}_h={0.01, 0.1, 1.0, 5.0}"

# --- make room for a new "noLT" / block_diag40 cauchy experiment column
#     under the "Latent Continous State Block Diagonal" section by
#     inserting 4 rows before the "Music Chord 1" section ---
$ws.Rows("29:32").Insert()

# --- fill in the new experiment notes in column L ---
$ws.Range("L25").Value = "cauchy"
$ws.Range("L26").Value = "script="
$ws.Range("L27").Value = "data_source=block_diag40_s2"
$ws.Range("L28").Value = "LT only"
$ws.Range("L29").Value = "lambda=1.6"
$ws.Range("L30").Value = "reps=5"
$ws.Range("L31").Value = "results/continuous_latent_syn/block_diag40/block_diag40_s2"
$ws.Range("L32").Value = "LT_hdp_hmm_w0_lambda1_cauchy"

# --- widen column L so the longer notes are readable ---
$ws.Columns("L").ColumnWidth = 52.33203125

# --- move the view back to the top-left and point the active selection at
#     the new notes column ---
$ws.Range("L24").Select()
